# issue #5: property land done
# Clean up stray spaces (and one stray hyphen) that crept into several
# shared-string cell values across the property-report sheets.

$wb = $excel.ActiveWorkbook

# 汽車 (Car) sheet
$ws1 = $wb.Worksheets.Item("汽車")
$ws1.Range("E2").Value = "96年07月01日"

# 存款 (Deposit) sheet
$ws2 = $wb.Worksheets.Item("存款")
$ws2.Range("B4").Value = "中國信託商業銀行城中分行"
$ws2.Range("B5").Value = "中國信託商業銀行城中分行"
$ws2.Range("B6").Value = "台新國際商業銀行鳳山分行"
$ws2.Range("B7").Value = "台新國際商業銀行鳳山分行"
$ws2.Range("B8").Value = "台新國際商業銀行鳳山分行"

# 基金受益憑證 (Fund) sheet
$ws4 = $wb.Worksheets.Item("基金受益憑證")
$ws4.Range("B2").Value = "貝萊德美國政府房貸債券基金A3美元（配現）"
$ws4.Range("B3").Value = "富蘭克林坦伯頓全球投資糸列一美國政府"
$ws4.Range("B4").Value = "摩根美國複合收益A股(人息)美元"
$ws4.Range("D2").Value = "中國信託商業銀行"
$ws4.Range("D3").Value = "中國信託商業銀行"
$ws4.Range("D4").Value = "中國信託商業銀行"

# 保險 (Insurance) sheet
$ws5 = $wb.Worksheets.Item("保險")
$ws5.Range("C2").Value = "新光人壽全心全意終身還本保險"
$ws5.Range("C3").Value = "世紀領航萬能終身壽險計劃A"
